# Apply edits described by the diff:
# 1. Update Summary sheet (B3:B9) with new aggregate figures after trade #26 closes.
# 2. Update Strategy Status sheet (C4:G4) MarketMaking row with new aggregate figures.
# 3. Append new trade #26 (row 27) to both "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199
$summary.Range("B4").Value = -1
$summary.Range("B5").Value = -0.77
$summary.Range("B6").Value = 26
$summary.Range("B8").Value = 15
$summary.Range("B9").Value = 23.08

# ---- Strategy Status sheet ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99
$status.Range("D4").Value = 26
$status.Range("E4").Value = -1
$status.Range("F4").Value = -1
$status.Range("G4").Value = 23.08

# ---- New trade row data ----
$tradeNum = 26
$date = "2026-02-17"
$time = "08:22:36"
$strategy = "MarketMaking"
$side = "UP"
$entryPrice = 0.66
$exitPrice = 0.561698
$status2 = "CLOSED"
$pnlPct = -14.8942
$pnlDollar = -0.1
$capitalAfter = 99
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.14

function Add-TradeRow($ws) {
    $ws.Range("A27").Value = $tradeNum
    $ws.Range("B27").Value = "'" + $date
    $ws.Range("C27").Value = $time
    $ws.Range("D27").Value = $strategy
    $ws.Range("E27").Value = $side
    $ws.Range("F27").Value = $entryPrice
    $ws.Range("G27").Value = $exitPrice
    $ws.Range("H27").Value = $status2
    $ws.Range("I27").Value = $pnlPct
    $ws.Range("J27").Value = $pnlDollar
    $ws.Range("K27").Value = $capitalAfter
    $ws.Range("L27").Value = $entrySlippage
    $ws.Range("M27").Value = $exitSlippage
    $ws.Range("N27").Value = $confidence
    $ws.Range("O27").Value = $entryReason
    $ws.Range("P27").Value = $exitReason
    $ws.Range("Q27").Value = $duration
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
